$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 BaggingClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                            colsample_bytree=0.5,
                                                            learning_rate=0.01,
                                                            max_depth=1,
                                                            min_child_samples=5,
                                                            num_leaves=2,
                                                            random_state=42,
                                                            subsample=0.7),
                                   random_state=42))])'
$ws.Range("B2").Value = 0.6984432234432234
$ws.Range("C2").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__subsample'': 0.7, ''model__estimator__num_leaves'': 2, ''model__estimator__min_child_samples'': 5, ''model__estimator__max_depth'': 1, ''model__estimator__learning_rate'': 0.01, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': None, ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D2").Value = 0.8552244409839967
$ws.Range("E2").Value = 0.5655480130980131
$ws.Range("F2").Value = 0.7777777777777777
$ws.Range("G2").Value = 0.859119203699992
$ws.Range("H2").Value = 0.5678051587301587
$ws.Range("I2").Value = 0.7
$ws.Range("J2").Value = 0.8605957446808512
$ws.Range("K2").Value = 0.5943333333333333
$ws.Range("L2").Value = 0.875
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 0 0 1 1 1 1 1 1 0 1 0 1 1 1 1 1 1 1]'
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f766340>),
                (''model'',
                 BaggingClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                            colsample_bytree=0.5,
                                                            learning_rate=0.05,
                                                            max_depth=5,
                                                            num_leaves=2,
                                                            random_state=42,
                                                            subsample=0.7),
                                   n_estimators=5, random_state=42))])'
$ws.Range("B3").Value = 0.6833699633699634
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f73c250>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__subsample'': 0.7, ''model__estimator__num_leaves'': 2, ''model__estimator__min_child_samples'': 20, ''model__estimator__max_depth'': 5, ''model__estimator__learning_rate'': 0.05, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': None, ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D3").Value = 0.8672354309505026
$ws.Range("E3").Value = 0.5362722277722277
$ws.Range("F3").Value = 0.8108108108108109
$ws.Range("G3").Value = 0.8573069239876172
$ws.Range("H3").Value = 0.597225
$ws.Range("I3").Value = 0.7142857142857143
$ws.Range("J3").Value = 0.8880212765957447
$ws.Range("K3").Value = 0.5208333333333333
$ws.Range("L3").Value = 0.9375
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 0 1 1 1 1 1 0 1 0 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f73c6a0>),
                (''model'',
                 BaggingClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                            colsample_bytree=0.5,
                                                            learning_rate=0.01,
                                                            max_depth=5,
                                                            min_child_samples=5,
                                                            num_leaves=2,
                                                            random_state=42,
                                                            subsample=0.5),
                                   random_state=42))])'
$ws.Range("B4").Value = 0.7077222777222777
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f510bb0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__subsample'': 0.5, ''model__estimator__num_leaves'': 2, ''model__estimator__min_child_samples'': 5, ''model__estimator__max_depth'': 5, ''model__estimator__learning_rate'': 0.01, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': None, ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D4").Value = 0.8798226901523074
$ws.Range("E4").Value = 0.5715195415695414
$ws.Range("F4").Value = 0.8717948717948718
$ws.Range("G4").Value = 0.8682606778324597
$ws.Range("H4").Value = 0.5916146825396826
$ws.Range("I4").Value = 0.85
$ws.Range("J4").Value = 0.8966444444444445
$ws.Range("K4").Value = 0.5876
$ws.Range("L4").Value = 0.8947368421052632
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 1 1 1 0 1 1 1 1 0 1 1 1 1 1 1 1 1 1 0 1 1]'
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 BaggingClassifier(estimator=LGBMClassifier(colsample_bytree=0.9,
                                                            learning_rate=0.01,
                                                            max_depth=1,
                                                            min_child_samples=10,
                                                            num_leaves=20,
                                                            random_state=42,
                                                            subsample=0.9),
                                   random_state=42))])'
$ws.Range("B5").Value = 0.7499999999999999
$ws.Range("C5").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__subsample'': 0.9, ''model__estimator__num_leaves'': 20, ''model__estimator__min_child_samples'': 10, ''model__estimator__max_depth'': 1, ''model__estimator__learning_rate'': 0.01, ''model__estimator__colsample_bytree'': 0.9, ''model__estimator__class_weight'': None, ''model__estimator__boosting_type'': ''gbdt''}'
$ws.Range("D5").Value = 0.8772628911239615
$ws.Range("E5").Value = 0.6095717782217782
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.8523355697401049
$ws.Range("H5").Value = 0.5867119047619047
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.9151020408163265
$ws.Range("K5").Value = 0.6576666666666666
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("O5").Value = 99

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 BaggingClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                            colsample_bytree=0.5,
                                                            learning_rate=0.01,
                                                            max_depth=1,
                                                            min_child_samples=5,
                                                            num_leaves=20,
                                                            random_state=42,
                                                            subsample=0.7),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B6").Value = 0.7288095238095238
$ws.Range("C6").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__subsample'': 0.7, ''model__estimator__num_leaves'': 20, ''model__estimator__min_child_samples'': 5, ''model__estimator__max_depth'': 1, ''model__estimator__learning_rate'': 0.01, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': None, ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D6").Value = 0.8872785024665917
$ws.Range("E6").Value = 0.6216035520035519
$ws.Range("F6").Value = 0.6285714285714286
$ws.Range("G6").Value = 0.8939285151979025
$ws.Range("H6").Value = 0.6065603174603175
$ws.Range("I6").Value = 0.4583333333333333
$ws.Range("J6").Value = 0.8907692307692308
$ws.Range("K6").Value = 0.6628333333333335
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("O6").Value = 89
